$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header row (id_tindakan / spesialisasi / nama_tindakan / tarif / keterangan)
# was removed - select it first (mirrors doing it by hand) then delete it,
# which shifts every data row up by one.
$ws.Rows(1).Select()
$ws.Rows(1).Delete()

# Fill in the newly added "id_tindakan" code column (A) for the three
# remaining data rows, and refresh a couple of descriptive texts in
# column E to their corrected/expanded wording.
$ws.Range("E1").Value = "Pemeriksaan biasa untuk menghasilkan diagnosa"
$ws.Range("E3").Value = "Pemberian imunisasi campak, difteri, tuberkulosis, dan polio"

$ws.Range("A1").Value = "T0001"
$ws.Range("A2").Value = "T0002"
$ws.Range("A3").Value = "T0003"

# The last two old data rows (Pemberian resep / Antibiotik) are no longer
# part of the table - clear them back to just the formatted-but-empty D
# cell that remains on rows 4 and 5.
$ws.Range("B4:E4").ClearContents()
$ws.Range("B5:E5").ClearContents()

# Leave the full former header row selected, same as after performing the
# row deletion interactively.
$ws.Range("A1:XFD1").Select()
